$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serials, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44308, 1, 9, 59.38242280285036),
    @(44309, 2, 10, 65.98046978094484),
    @(44310, 0, 9, 59.38242280285036),
    @(44311, 3, 8, 52.78437582475588),
    @(44312, 1, 9, 59.38242280285036)
)

$startRow = 234
$endRow = $startRow + $data.Length - 1

# Apply the same formatting used on the rest of column A (date number
# format, bold font, thin border, centered alignment) to the new cells
# before filling in values.
$ws.Range("A233").Copy($ws.Range("A${startRow}:A${endRow}"))

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
